$p = $ppt.ActivePresentation

function Get-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $candidate = $slide.Shapes.Item($i)
        if ($candidate.Name -eq $name) {
            return $candidate
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# Slide 5 ("Timeline & Milestones") - Table Placeholder 3
#   - Header row (row 1, sz=1400): remove bold from all 4 cells
#   - "Phase 3" row (row 4, sz=1100): remove bold from all 4 cells
# ---------------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$table5 = (Get-ShapeByName $slide5 "Table Placeholder 3").Table

for ($c = 1; $c -le $table5.Columns.Count; $c++) {
    $table5.Cell(1, $c).Shape.TextFrame.TextRange.Font.Bold = 0
}
for ($c = 1; $c -le $table5.Columns.Count; $c++) {
    $table5.Cell(4, $c).Shape.TextFrame.TextRange.Font.Bold = 0
}

# ---------------------------------------------------------------------------
# Slide 8 ("Investment Summary") - Table Placeholder 3
#   - Header row (row 1, sz=1400): remove bold from all 7 cells
#   - (TOTAL INVESTMENT row stays bold - not touched)
# ---------------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)
$table8 = (Get-ShapeByName $slide8 "Table Placeholder 3").Table

for ($c = 1; $c -le $table8.Columns.Count; $c++) {
    $table8.Cell(1, $c).Shape.TextFrame.TextRange.Font.Bold = 0
}

# ---------------------------------------------------------------------------
# Slide 9 ("Next Steps") - Text Placeholder 3
#   - Each paragraph "Label: rest of sentence" becomes two runs:
#     a bold "Label:" run followed by a regular " rest of sentence" run.
# ---------------------------------------------------------------------------
$slide9 = $p.Slides.Item(9)
$tf9 = (Get-ShapeByName $slide9 "Text Placeholder 3").TextFrame
$tr9 = $tf9.TextRange

$labels = @("Decision:", "Kickoff:", "Team Formation:", "Week 1:", "Week 2:", "Week 3:", "Week 4:")

for ($i = 1; $i -le $labels.Count; $i++) {
    $para = $tr9.Paragraphs($i)
    $label = $labels[$i - 1]
    $boldPart = $para.Characters(1, $label.Length)
    $boldPart.Font.Bold = 1
}
